$wb = $excel.ActiveWorkbook

# Mapping of row -> new value for column F ("想去人数")
$updates = @{
    5  = 96
    6  = 133
    7  = 1290
    8  = 1547
    10 = 408
    12 = 165
    16 = 272
    17 = 311
    19 = 1747
    23 = 675
    26 = 4212
    29 = 1102
    30 = 491
    32 = 573
    34 = 274
}

# The update applies identically to sheets "展览" (1st) and "全部类型" (4th)
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
